# Apply OCR-autocorrect style text fixes to the shared strings used
# across the "Simple Fields" / "Simple Fields - Formatted" and
# "Items" / "Items - Formatted" worksheets.

$wb = $excel.ActiveWorkbook

# --- Vendor Address / Currency fixes (Simple Fields + Simple Fields - Formatted) ---
$simpleSheets = @("Simple Fields", "Simple Fields - Formatted")
foreach ($sheetName in $simpleSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("B2").Value = "1580 NW Gilman Blvd Suite 1 Issaquah WA, 98027"
    $ws.Range("I2").Value = "SGD"
}

# --- Item description fixes (Items + Items - Formatted) ---
$itemSheets = @("Items", "Items - Formatted")
foreach ($sheetName in $itemSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A2").Value  = "green onion Pancakes ÂY/MAf (1)"
    $ws.Range("A3").Value  = "Pan Fried Leek Dumplings #7 (2)"
    $ws.Range("A4").Value  = "Pork Xiao Long Bao(10) AP])\`$E(10)"
    $ws.Range("A5").Value  = "Q-BA( (5) ĦEH'L (5)"
    $ws.Range("A6").Value  = "Chicken potstickers HÈP]`$9I5(6)"
    $ws.Range("A7").Value  = "Tomato Mushroom Steamed dumpli PEÅINABEEMKK (6)"
    $ws.Range("A8").Value  = "Zucchini shrimp dumplings ĦJU]K"
    $ws.Range("A9").Value  = 'beef stew nodle soup (Non Spicy "H751PJ(74k)'
    $ws.Range("A10").Value = "dandan noodle"
    $ws.Range("A11").Value = "banana naan bread ¥"
    $ws.Range("A12").Value = "house made plum juice ĚUNNT"
}
